$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.717.98"
$ws.Range("E2").Value = "  -0.59%  "

# Row 3
$ws.Range("D3").Value = "2.212.64"
$ws.Range("E3").Value = "  -1.13%  "

# Row 4
$ws.Range("E4").Value = "  -0.27%  "

# Row 5
$ws.Range("D5").Value = "'253.34"
$ws.Range("E5").Value = "  +2.56%  "

# Row 6
$ws.Range("E6").Value = "  -0.44%  "

# Row 7
$ws.Range("D7").Value = "'75.49"
$ws.Range("E7").Value = "  -0.71%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("E9").Value = "  -4.33%  "

# Row 10
$ws.Range("D10").Value = "'40.91"
$ws.Range("E10").Value = "  +0.18%  "

# Row 11
$ws.Range("D11").Value = "'0.0921"
$ws.Range("E11").Value = "  -1.12%  "

# Row 12
$ws.Range("E12").Value = "  -1.03%  "

# Row 13
$ws.Range("E13").Value = "  +0.32%  "

# Row 14
$ws.Range("D14").Value = "2.534.02"
$ws.Range("E14").Value = "  -0.70%  "

# Row 15
$ws.Range("D15").Value = "'14.33"
$ws.Range("E15").Value = "  -2.45%  "

# Row 16
$ws.Range("D16").Value = "2.207.72"
$ws.Range("E16").Value = "  -1.37%  "

# Row 17
$ws.Range("D17").Value = "'0.782"
$ws.Range("E17").Value = "  -3.80%  "

# Row 18
$ws.Range("D18").Value = "42.638.68"
$ws.Range("E18").Value = "  -0.53%  "

# Row 19
$ws.Range("E19").Value = "  -1.71%  "

# Row 20
$ws.Range("D20").Value = "'71.13"
$ws.Range("E20").Value = "  +0.02%  "

# Row 21
$ws.Range("D21").Value = "'5.97"
$ws.Range("E21").Value = "  -0.38%  "

# Row 22
$ws.Range("D22").Value = "'2.20"
$ws.Range("E22").Value = "  -0.27%  "

# Row 23
$ws.Range("D23").Value = "'228.89"
$ws.Range("E23").Value = "  -0.59%  "

# Row 24
$ws.Range("D24").Value = "'9.43"
$ws.Range("E24").Value = "  -9.27%  "

# Row 25
$ws.Range("E25").Value = "  +0.03%  "

# Row 26
$ws.Range("D26").Value = "'10.58"
$ws.Range("E26").Value = "  -2.99%  "

# Row 27
$ws.Range("D27").Value = "'3.36"
$ws.Range("E27").Value = "  -0.16%  "

# Row 28
$ws.Range("D28").Value = "'39.21"
$ws.Range("E28").Value = "  +4.18%  "

# Row 29
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'2.18"
$ws.Range("E29").Value = "  -3.34%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.19"
$ws.Range("E30").Value = "  -0.89%  "

# Row 31
$ws.Range("D31").Value = "'173.43"
$ws.Range("E31").Value = "  +0.00%  "

# Row 32
$ws.Range("D32").Value = "'20.23"
$ws.Range("E32").Value = "  -0.31%  "

# Row 33
$ws.Range("D33").Value = "'0.0830"
$ws.Range("E33").Value = "  +4.55%  "

# Row 34
$ws.Range("E34").Value = "  -3.07%  "

# Row 35
$ws.Range("E35").Value = "  -1.10%  "

# Row 36
$ws.Range("D36").Value = "'0.110"
$ws.Range("E36").Value = "  -1.87%  "

# Row 37
$ws.Range("E37").Value = "  +4.97%  "

# Row 38
$ws.Range("D38").Value = "'4.28"
$ws.Range("E38").Value = "  -1.47%  "

# Row 39
$ws.Range("D39").Value = "'12.32"
$ws.Range("E39").Value = "  -5.37%  "

# Row 40
$ws.Range("D40").Value = "'2.09"
$ws.Range("E40").Value = "  -2.16%  "

# Row 41
$ws.Range("E41").Value = "  +16.92%  "

# Row 42
$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D42").Value = "'5.25"
$ws.Range("E42").Value = "  -5.72%  "

# Row 43
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").Value = "'59.90"
$ws.Range("E43").Value = "  -0.19%  "

# Row 44
$ws.Range("E44").Value = "  -1.65%  "

# Row 45
$ws.Range("D45").Value = "'102.86"
$ws.Range("E45").Value = "  -2.36%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'8.37"
$ws.Range("E46").Value = "  -2.62%  "

# Row 47
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.0982"
$ws.Range("E47").Value = "  -0.79%  "

# Row 48
$ws.Range("D48").Value = "'0.459"
$ws.Range("E48").Value = "  +3.80%  "

# Row 49
$ws.Range("E49").Value = "  -0.28%  "

# Row 50
$ws.Range("E50").Value = "  -0.89%  "

# Row 51
$ws.Range("E51").Value = "  -0.93%  "

Write-Host "Applied cryptos update"
